$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Rows 12 and 13 currently hold:
#   Row 12: "Enterprises (absolute #)"             / "2646549"
#   Row 13: "Enterprises density (per 1000 people)" / "39.8"
# Swap them so the density row comes first (matches the new
# sharedStrings ordering in the target workbook).
#
# "2646549" and "39.8" are stored as text, not numbers. A bare
# Value assignment of a numeric-looking string gets auto-coerced to a
# number, so prefix with an apostrophe to force text entry, then reset
# the cell Style afterwards so the quote-prefix formatting doesn't stick
# around on the cell (keeps the original plain-text/"Normal" styling).

$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("C12").Value = "'39.8"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'39.8"
$ws.Range("D12").Style = "Normal"

$ws.Range("A13").Value = "Enterprises (absolute #)"
$ws.Range("C13").Value = "'2646549"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'2646549"
$ws.Range("D13").Style = "Normal"
